# Apply the edit described by the diff: insert 6 new region rows into the
# "10 March 2023" sheet, shifting the existing rows below each insertion
# point down by one. We insert from the bottom of the sheet upward so that
# each insertion position (expressed in terms of the *original*, pre-edit
# row numbers) stays valid for the remaining insertions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-DataRow {
    param(
        [int]$RowNumber,
        [string]$Region,
        [string]$Issues,
        [int]$Minor,
        [string]$MinorTimestamp,
        [int]$Major,
        [string]$MajorTimestamp,
        [bool]$NativeEmbassies,
        [string]$Link,
        [string]$Organizations
    )

    $ws.Rows.Item($RowNumber).Insert()

    $a = $ws.Cells.Item($RowNumber, 1)
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Borders.LineStyle = 1

    $a.Value = $Region
    $ws.Cells.Item($RowNumber, 2).Value = $Issues
    $ws.Cells.Item($RowNumber, 3).Value = $Minor
    $ws.Cells.Item($RowNumber, 4).Value = $MinorTimestamp
    $ws.Cells.Item($RowNumber, 5).Value = $Major
    $ws.Cells.Item($RowNumber, 6).Value = $MajorTimestamp
    $ws.Cells.Item($RowNumber, 7).Value = $NativeEmbassies
    $ws.Cells.Item($RowNumber, 8).Value = $Link
    $ws.Cells.Item($RowNumber, 9).Value = $Organizations
}

# Insertions, processed from the bottom of the sheet upward. The row
# number passed in is "immediately after" the last unmodified row at the
# time of insertion (i.e. original_row + 1), which is why the numbers
# below are not monotonically shifted for later edits.

# After original row 143 ("Ben") -> new row 144: The Fifth Sovereign Charter
Insert-DataRow 144 "The Fifth Sovereign Charter" "WFE, RO" 3196 "0:53:16" 4794 "1:19:54" $false "https://www.nationstates.net/region=the_fifth_sovereign_charter" "Unknown"

# After original row 136 ("The Brotherhood of Malice Girl Scouts") -> new row 137: raiding am I right
Insert-DataRow 137 "raiding am I right" "RO" 2278 "0:37:58" 3417 "0:56:57" $false "https://www.nationstates.net/region=raiding_am_i_right" "Unknown"

# After original row 129 ("Land of Places") -> new row 130: Nowhere
Insert-DataRow 130 "Nowhere" "WFE, RO" 1341 "0:22:21" 2012 "0:33:32" $false "https://www.nationstates.net/region=nowhere" "Unknown"

# After original row 114 ("Cyberius Confederation") -> new row 115: Coffee House
Insert-DataRow 115 "Coffee House" "WFE, RO" 1091 "0:18:11" 1636 "0:27:16" $false "https://www.nationstates.net/region=coffee_house" "Unknown"

# After original row 107 ("The 8th Grade Isles") -> new row 108: Open Ocean 9
Insert-DataRow 108 "Open Ocean 9" "WFE" 775 "0:12:55" 1163 "0:19:23" $false "https://www.nationstates.net/region=open_ocean_9" "Unknown"

# After original row 74 ("Archination") -> new row 75: NPS
Insert-DataRow 75 "NPS" "WFE" 198 "0:03:18" 296 "0:04:56" $false "https://www.nationstates.net/region=nps" "Unknown"
